$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet / tenant
$ws.Name = "Partu"

# Update property address
$ws.Range("A2").Value = "Heinrichstraße  46 , 59077  Hamm,Reihenmittelhaus,"

# Update tenant name
$ws.Range("B3").Value = "Partu"

# Row 5: Entwässerung Fläche -> Entwässerung, value 168,96 -> 84,00
$ws.Range("A5").Value = "Entwässerung"
$ws.Range("C5").Value = "84,00"
$ws.Range("F5").Value = "84,00"

# Row 6: Grundsteuer 318,85 -> 301,50
$ws.Range("A6").Value = "Grundsteuer"
$ws.Range("C6").Value = "301,50"
$ws.Range("F6").Value = "301,50"

# Row 7: Müllabfuhr 118,20 -> 147,62
$ws.Range("A7").Value = "Müllabfuhr"
$ws.Range("C7").Value = "147,62"
$ws.Range("F7").Value = "147,62"

# Row 8: Gebäudeversicherung 280,00 -> 231,95
$ws.Range("A8").Value = "Gebäudeversicherung"
$ws.Range("C8").Value = "231,95"
$ws.Range("F8").Value = "231,95"

# Row 9: Haftpflichtversicherung 43,91 -> 46,10
$ws.Range("A9").Value = "Haftpflichtversicherung"
$ws.Range("C9").Value = "46,10"
$ws.Range("F9").Value = "46,10"

# Row 10: Heizungswartung 121,46 -> 163,60
$ws.Range("A10").Value = "Heizungswartung"
$ws.Range("C10").Value = "163,60"
$ws.Range("F10").Value = "163,60"

# Row 11: new Kabelfernsehen 146,28 (B/D/E already hold "Haushalt "/"1"/"1")
$ws.Range("A11").Value = "Kabelfernsehen"
$ws.Range("C11").Value = "146,28"
$ws.Range("F11").Value = "146,28"

# Row 12: Rauchmelder 30,37 -> 111,45
$ws.Range("A12").Value = "Rauchmelder"
$ws.Range("C12").Value = "111,45"
$ws.Range("F12").Value = "111,45"

# Row 13: new Schornsteinfeger 64,55
$ws.Range("A13").Value = "Schornsteinfeger"
$ws.Range("C13").Value = "64,55"
$ws.Range("F13").Value = "64,55"

# Row 14: Summe total 1.919,47 -> 1.297,05
$ws.Range("F14").Value = "1.297,05"
